$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.756.00'
$ws.Range('E2').Value = '  -2.81%  '
$ws.Range('D3').Value = '2.575.96'
$ws.Range('E3').Value = '  -4.97%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '546.82'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '154.18'
$ws.Range('E6').Value = '  -1.68%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  +1.60%  '
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '5.38'
$ws.Range('E11').Value = '  -2.26%  '
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').Value = '3.032.50'
$ws.Range('E13').Value = '  -5.20%  '
$ws.Range('D14').Value = '25.38'
$ws.Range('E14').Value = '  -4.19%  '
$ws.Range('D15').Value = '61.670.45'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('E17').Value = '  -5.29%  '
$ws.Range('D18').Value = '11.51'
$ws.Range('E18').Value = '  -5.23%  '
$ws.Range('D19').Value = '4.53'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').Value = '336.59'
$ws.Range('E20').Value = '  -3.57%  '
$ws.Range('D21').Value = '6.05'
$ws.Range('E21').Value = '  -5.44%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '0.492'
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('D24').Value = '63.22'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '8.07'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = '7.50'
$ws.Range('E28').Value = '  +5.44%  '
$ws.Range('D29').Value = '0.0₃0836'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('E30').Value = '  -2.33%  '
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('D32').Value = '160.52'
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '4.73'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').Value = '19.15'
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '332.56'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('D39').Value = '0.921'
$ws.Range('E39').Value = '  -3.07%  '
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('D41').Value = '3.94'
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').Value = '2.127.34'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('D48').Value = '19.56'
$ws.Range('E48').Value = '  -4.65%  '
$ws.Range('D49').Value = '0.0544'
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('D50').Value = '0.0964'
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('E51').Value = '  -2.25%  '
